# The commit swaps the contents of ppt/theme/theme1.xml ("Office Theme")
# and ppt/theme/theme2.xml ("Integral") — i.e. the deck's editable/active
# theme (the one behind the slide master, reachable through the
# PowerPoint object model as the presentation's ThemeColorScheme) changes
# its 12 scheme colors from the "Integral" palette to the "Office Theme"
# palette.
#
# Colors are addressed in DrawingML clrScheme order via
# Slide.ThemeColorScheme, which edits the 12 colors (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) in place without disturbing any other part
# of the theme part (font scheme / format scheme / names stay intact).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target values = the "Office Theme" color scheme (formerly theme1.xml),
# expressed as VBA-style RGB() integers (0xBBGGRR).
$tcs.Item(1).RGB  = 0        # dk1      000000
$tcs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388  # dk2      44546A
$tcs.Item(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407    # accent4  FFC000
$tcs.Item(9).RGB  = 12874308 # accent5  4472C4
$tcs.Item(10).RGB = 4697456  # accent6  70AD47
$tcs.Item(11).RGB = 12673797 # hlink    0563C1
$tcs.Item(12).RGB = 7491477  # folHlink 954F72
